$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E retain text formatting so numeric-looking
# strings (e.g. "1.004") are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.989.31"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "1.910.07"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").Value = "0.7775"
$ws.Range("E5").Value = "  +4.47%  "
$ws.Range("D6").Value = "241.85"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").Value = "0.3143"
$ws.Range("E8").Value = "  +2.36%  "
$ws.Range("D9").Value = "25.98"
$ws.Range("E9").Value = "  +1.09%  "
$ws.Range("D10").Value = "0.06878"
$ws.Range("E10").Value = "  -0.34%  "
$ws.Range("D11").Value = "0.07959"
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("D12").Value = "1.899.29"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").Value = "0.7396"
$ws.Range("E13").Value = "  -2.28%  "
$ws.Range("D14").Value = "5.195"
$ws.Range("E14").Value = "  -0.79%  "
$ws.Range("D15").Value = "92.75"
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("D16").Value = "30.006.30"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "13.91"
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").Value = "5.859"
$ws.Range("E18").Value = "  -4.68%  "
$ws.Range("D19").Value = "245.83"
$ws.Range("E19").Value = "  +3.86%  "
$ws.Range("D20").Value = "0.000007731"
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").Value = "2.149.97"
$ws.Range("E22").Value = "  -0.59%  "
$ws.Range("D23").Value = "1.004"
$ws.Range("E23").Value = "  +0.37%  "
$ws.Range("D24").Value = "6.856"
$ws.Range("D25").Value = "168.66"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").Value = "9.262"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("D27").Value = "0.1364"
$ws.Range("E27").Value = "  +7.71%  "
$ws.Range("D28").Value = "18.90"
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("D29").Value = "2.023"
$ws.Range("E29").Value = "  -1.56%  "
$ws.Range("D30").Value = "1.375"
$ws.Range("E30").Value = "  +1.72%  "
$ws.Range("D31").Value = "1.518"
$ws.Range("D32").Value = "4.312"
$ws.Range("E32").Value = "  +0.41%  "
$ws.Range("D33").Value = "4.071"
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("D34").Value = "0.05448"
$ws.Range("E34").Value = "  +2.36%  "
$ws.Range("D35").Value = "1.251"
$ws.Range("E35").Value = "  -2.85%  "
$ws.Range("D36").Value = "0.7316"
$ws.Range("E36").Value = "  -1.18%  "
$ws.Range("D37").Value = "2.734"
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("D38").Value = "0.01930"
$ws.Range("E38").Value = "  -0.68%  "
$ws.Range("D39").Value = "2.797"
$ws.Range("E39").Value = "  +1.21%  "
$ws.Range("D40").Value = "6.119"
$ws.Range("E40").Value = "  -2.07%  "
$ws.Range("D41").Value = "0.4401"
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("D42").Value = "71.84"
$ws.Range("E42").Value = "  -1.10%  "
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("D44").Value = "0.8378"
$ws.Range("E44").Value = "  +0.80%  "
$ws.Range("D45").Value = "1.867"
$ws.Range("E45").Value = "  -4.57%  "
$ws.Range("D46").Value = "100.24"
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "9.750"
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "7.498"
$ws.Range("E48").Value = "  -2.69%  "
$ws.Range("D49").Value = "978.79"
$ws.Range("E49").Value = "  +8.37%  "
$ws.Range("D50").Value = "2.056.83"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").Value = "36.12"
$ws.Range("E51").Value = "  -1.31%  "
